# ===========================================================================
# river update May 2024
# Refreshes the "Manawatu at Whirokino" state-results sheet:
#   - Recomputed Mean ("G") values (and a few Median/95th-percentile/etc.
#     values) across the 2019-2023 rolling-window blocks using newer source data
#   - Replaces the placeholder 2018-2022 E.coli block (rows 109-112, 124-127)
#     with corrected figures
#   - Appends the new 2019-2023 reporting block (rows 385-401)
# ===========================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update existing cells with refreshed statistics ---
$ws.Cells.Item(27, 7).Value = 0.0419002617283696
$ws.Cells.Item(28, 7).Value = 0.0419002617283696
$ws.Cells.Item(36, 7).Value = 0.0435332376528311
$ws.Cells.Item(37, 7).Value = 0.0435332376528311
$ws.Cells.Item(45, 7).Value = 0.0477795504037192
$ws.Cells.Item(46, 7).Value = 0.0477795504037192
$ws.Cells.Item(54, 7).Value = 0.0494262206261609
$ws.Cells.Item(55, 7).Value = 0.0494262206261609
$ws.Cells.Item(65, 7).Value = 0.0498332163720666
$ws.Cells.Item(66, 7).Value = 0.0498332163720666
$ws.Cells.Item(106, 7).Value = 0.367721786559556
$ws.Cells.Item(109, 7).Value = 1350.83166172185
$ws.Cells.Item(109, 9).Value = 5365.0
$ws.Cells.Item(109, 14).Value = 4490.0
$ws.Cells.Item(110, 7).Value = 1350.83166172185
$ws.Cells.Item(110, 9).Value = 5365.0
$ws.Cells.Item(110, 14).Value = 4490.0
$ws.Cells.Item(111, 7).Value = 1350.83166172185
$ws.Cells.Item(111, 9).Value = 5365.0
$ws.Cells.Item(111, 14).Value = 4490.0
$ws.Cells.Item(112, 7).Value = 1350.83166172185
$ws.Cells.Item(112, 9).Value = 5365.0
$ws.Cells.Item(112, 14).Value = 4490.0
$ws.Cells.Item(121, 7).Value = 0.346251330092561
$ws.Cells.Item(124, 7).Value = 1345.51780069316
$ws.Cells.Item(124, 9).Value = 5500.0
$ws.Cells.Item(124, 14).Value = 5020.0
$ws.Cells.Item(125, 7).Value = 1345.51780069316
$ws.Cells.Item(125, 9).Value = 5500.0
$ws.Cells.Item(125, 14).Value = 5020.0
$ws.Cells.Item(126, 7).Value = 1345.51780069316
$ws.Cells.Item(126, 9).Value = 5500.0
$ws.Cells.Item(126, 14).Value = 5020.0
$ws.Cells.Item(127, 7).Value = 1345.51780069316
$ws.Cells.Item(127, 9).Value = 5500.0
$ws.Cells.Item(127, 14).Value = 5020.0
$ws.Cells.Item(128, 7).Value = 0.0321700633500733
$ws.Cells.Item(129, 7).Value = 0.0321700633500733
$ws.Cells.Item(136, 7).Value = 0.342178210362842
$ws.Cells.Item(143, 7).Value = 0.0320076235802917
$ws.Cells.Item(144, 7).Value = 0.0320076235802917
$ws.Cells.Item(147, 7).Value = 0.640203343998376
$ws.Cells.Item(148, 7).Value = 0.640203343998376
$ws.Cells.Item(151, 7).Value = 0.429856781791413
$ws.Cells.Item(152, 7).Value = 0.0321536655322381
$ws.Cells.Item(153, 7).Value = 0.0321536655322381
$ws.Cells.Item(158, 7).Value = 0.0311597385128532
$ws.Cells.Item(159, 7).Value = 0.0311597385128532
$ws.Cells.Item(162, 7).Value = 0.592916621598403
$ws.Cells.Item(163, 7).Value = 0.592916621598403
$ws.Cells.Item(166, 7).Value = 0.502892496077127
$ws.Cells.Item(167, 7).Value = 0.0297833886768523
$ws.Cells.Item(168, 7).Value = 0.0297833886768523
$ws.Cells.Item(173, 7).Value = 0.0283471783661455
$ws.Cells.Item(174, 7).Value = 0.0283471783661455
$ws.Cells.Item(177, 7).Value = 0.577881310100071
$ws.Cells.Item(178, 7).Value = 0.577881310100071
$ws.Cells.Item(181, 7).Value = 0.504071337381382
$ws.Cells.Item(182, 6).Value = 0.02322
$ws.Cells.Item(182, 7).Value = 0.0263299011925891
$ws.Cells.Item(182, 12).Value = 0.02561
$ws.Cells.Item(183, 6).Value = 0.02322
$ws.Cells.Item(183, 7).Value = 0.0263299011925891
$ws.Cells.Item(183, 12).Value = 0.02561
$ws.Cells.Item(188, 7).Value = 0.0356552176199147
$ws.Cells.Item(189, 7).Value = 0.0356552176199147
$ws.Cells.Item(192, 7).Value = 0.586796564337359
$ws.Cells.Item(193, 7).Value = 0.586796564337359
$ws.Cells.Item(198, 7).Value = 0.578452393071266
$ws.Cells.Item(199, 7).Value = 0.0254993927180128
$ws.Cells.Item(199, 12).Value = 0.02561
$ws.Cells.Item(200, 7).Value = 0.0254993927180128
$ws.Cells.Item(200, 12).Value = 0.02561
$ws.Cells.Item(205, 7).Value = 0.046094627715921
$ws.Cells.Item(206, 7).Value = 0.046094627715921
$ws.Cells.Item(209, 7).Value = 0.615440632133969
$ws.Cells.Item(210, 7).Value = 0.615440632133969
$ws.Cells.Item(216, 7).Value = 0.0224824435654704
$ws.Cells.Item(217, 7).Value = 0.0224824435654704
$ws.Cells.Item(222, 7).Value = 0.0441851156893723
$ws.Cells.Item(223, 7).Value = 0.0441851156893723
$ws.Cells.Item(233, 7).Value = 0.0196312582784486
$ws.Cells.Item(234, 7).Value = 0.0196312582784486
$ws.Cells.Item(250, 7).Value = 0.0177207373071412
$ws.Cells.Item(251, 7).Value = 0.0177207373071412
$ws.Cells.Item(256, 7).Value = 0.0384751809519535
$ws.Cells.Item(256, 12).Value = 0.02283
$ws.Cells.Item(257, 7).Value = 0.0384751809519535
$ws.Cells.Item(257, 12).Value = 0.02283
$ws.Cells.Item(258, 7).Value = 0.502219922619336
$ws.Cells.Item(259, 7).Value = 0.502219922619336
$ws.Cells.Item(260, 7).Value = 0.583121991596377
$ws.Cells.Item(261, 7).Value = 0.583121991596377
$ws.Cells.Item(267, 7).Value = 0.0180448560597759
$ws.Cells.Item(268, 7).Value = 0.0180448560597759
$ws.Cells.Item(273, 6).Value = 0.02252
$ws.Cells.Item(273, 7).Value = 0.0323401653148127
$ws.Cells.Item(273, 12).Value = 0.02126
$ws.Cells.Item(274, 6).Value = 0.02252
$ws.Cells.Item(274, 7).Value = 0.0323401653148127
$ws.Cells.Item(274, 12).Value = 0.02126
$ws.Cells.Item(275, 7).Value = 0.52205920505532
$ws.Cells.Item(276, 7).Value = 0.52205920505532
$ws.Cells.Item(277, 7).Value = 0.599386163784538
$ws.Cells.Item(278, 7).Value = 0.599386163784538
$ws.Cells.Item(284, 7).Value = 0.0182448560597759
$ws.Cells.Item(285, 7).Value = 0.0182448560597759
$ws.Cells.Item(290, 6).Value = 0.02204
$ws.Cells.Item(290, 7).Value = 0.0320032421703398
$ws.Cells.Item(290, 12).Value = 0.02126
$ws.Cells.Item(290, 13).Value = 0.04837
$ws.Cells.Item(291, 6).Value = 0.02204
$ws.Cells.Item(291, 7).Value = 0.0320032421703398
$ws.Cells.Item(291, 12).Value = 0.02126
$ws.Cells.Item(291, 13).Value = 0.04837
$ws.Cells.Item(292, 7).Value = 0.544909205055319
$ws.Cells.Item(293, 7).Value = 0.544909205055319
$ws.Cells.Item(294, 7).Value = 0.613402830451205
$ws.Cells.Item(295, 7).Value = 0.613402830451205
$ws.Cells.Item(301, 7).Value = 0.0183115227264425
$ws.Cells.Item(302, 7).Value = 0.0183115227264425
$ws.Cells.Item(307, 7).Value = 0.0362679898493818
$ws.Cells.Item(308, 7).Value = 0.0362679898493818
$ws.Cells.Item(309, 7).Value = 0.551777693269125
$ws.Cells.Item(310, 7).Value = 0.551777693269125
$ws.Cells.Item(311, 7).Value = 0.610269497117872
$ws.Cells.Item(312, 7).Value = 0.610269497117872
$ws.Cells.Item(317, 7).Value = 0.369605326029721
$ws.Cells.Item(318, 7).Value = 0.0196448560597759
$ws.Cells.Item(319, 7).Value = 0.0196448560597759
$ws.Cells.Item(324, 7).Value = 0.0399554511982778
$ws.Cells.Item(325, 7).Value = 0.0399554511982778
$ws.Cells.Item(326, 7).Value = 0.564756026602459
$ws.Cells.Item(327, 7).Value = 0.564756026602459
$ws.Cells.Item(328, 7).Value = 0.628069497117872
$ws.Cells.Item(329, 7).Value = 0.628069497117872
$ws.Cells.Item(334, 7).Value = 0.359172999051949
$ws.Cells.Item(335, 7).Value = 0.0202705407514923
$ws.Cells.Item(336, 7).Value = 0.0202705407514923
$ws.Cells.Item(341, 7).Value = 0.0431573305578377
$ws.Cells.Item(342, 7).Value = 0.0431573305578377
$ws.Cells.Item(343, 7).Value = 0.572499417913576
$ws.Cells.Item(344, 7).Value = 0.572499417913576
$ws.Cells.Item(345, 7).Value = 0.64165848019135
$ws.Cells.Item(346, 7).Value = 0.64165848019135
$ws.Cells.Item(351, 7).Value = 0.376971263274444
$ws.Cells.Item(358, 7).Value = 0.0435281444575807
$ws.Cells.Item(359, 7).Value = 0.0435281444575807
$ws.Cells.Item(360, 7).Value = 0.527256303382954
$ws.Cells.Item(361, 7).Value = 0.527256303382954
$ws.Cells.Item(368, 7).Value = 0.364605326029721
$ws.Cells.Item(375, 7).Value = 0.0359745142355603
$ws.Cells.Item(376, 7).Value = 0.0359745142355603
$ws.Cells.Item(377, 7).Value = 0.550379110400497
$ws.Cells.Item(378, 7).Value = 0.550379110400497

# --- 2. Append the new 2019-2023 reporting rows (385-401) ---
# Row 385: Visual Clarity (Sediment class 3)
$ws.Cells.Item(385, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(385, 2).Value = "Visual Clarity (Sediment class 3)"
$ws.Cells.Item(385, 3).Value = "D"
$ws.Cells.Item(385, 4).Value = "2019 - 2023"
$ws.Cells.Item(385, 5).Value = "RepSite"
$ws.Cells.Item(385, 6).Value = 0.31
$ws.Cells.Item(385, 7).Value = 0.348421116736186
$ws.Cells.Item(385, 8).Value = 1.1
$ws.Cells.Item(385, 9).Value = 0.8975
$ws.Cells.Item(385, 10).Value = $null
$ws.Cells.Item(385, 11).Value = $null
$ws.Cells.Item(385, 12).Value = 0.34
$ws.Cells.Item(385, 13).Value = 0.5915
$ws.Cells.Item(385, 14).Value = 0.821
$ws.Cells.Item(385, 15).Value = 1792183
$ws.Cells.Item(385, 16).Value = 5512989
$ws.Cells.Item(385, 17).Value = "Horowhenua District"
$ws.Cells.Item(385, 18).Value = "Manawatū"
$ws.Cells.Item(385, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(385, 20).Value = "Mana_13a"
$ws.Cells.Item(385, 21).Value = "m"

# Row 386: DRP (95th Percentile)
$ws.Cells.Item(386, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(386, 2).Value = "DRP (95th Percentile)"
$ws.Cells.Item(386, 3).Value = "C"
$ws.Cells.Item(386, 4).Value = "2019 - 2023"
$ws.Cells.Item(386, 5).Value = "RepSite"
$ws.Cells.Item(386, 6).Value = 0.023
$ws.Cells.Item(386, 7).Value = 0.0217894736842105
$ws.Cells.Item(386, 8).Value = 0.05
$ws.Cells.Item(386, 9).Value = 0.03495
$ws.Cells.Item(386, 10).Value = $null
$ws.Cells.Item(386, 11).Value = $null
$ws.Cells.Item(386, 12).Value = 0.022
$ws.Cells.Item(386, 13).Value = 0.028
$ws.Cells.Item(386, 14).Value = 0.03088
$ws.Cells.Item(386, 15).Value = 1792183
$ws.Cells.Item(386, 16).Value = 5512989
$ws.Cells.Item(386, 17).Value = "Horowhenua District"
$ws.Cells.Item(386, 18).Value = "Manawatū"
$ws.Cells.Item(386, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(386, 20).Value = "Mana_13a"
$ws.Cells.Item(386, 21).Value = "mg/L"

# Row 387: DRP (Median)
$ws.Cells.Item(387, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(387, 2).Value = "DRP (Median)"
$ws.Cells.Item(387, 3).Value = "D"
$ws.Cells.Item(387, 4).Value = "2019 - 2023"
$ws.Cells.Item(387, 5).Value = "RepSite"
$ws.Cells.Item(387, 6).Value = 0.023
$ws.Cells.Item(387, 7).Value = 0.0217894736842105
$ws.Cells.Item(387, 8).Value = 0.05
$ws.Cells.Item(387, 9).Value = 0.03495
$ws.Cells.Item(387, 10).Value = $null
$ws.Cells.Item(387, 11).Value = $null
$ws.Cells.Item(387, 12).Value = 0.022
$ws.Cells.Item(387, 13).Value = 0.028
$ws.Cells.Item(387, 14).Value = 0.03088
$ws.Cells.Item(387, 15).Value = 1792183
$ws.Cells.Item(387, 16).Value = 5512989
$ws.Cells.Item(387, 17).Value = "Horowhenua District"
$ws.Cells.Item(387, 18).Value = "Manawatū"
$ws.Cells.Item(387, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(387, 20).Value = "Mana_13a"
$ws.Cells.Item(387, 21).Value = "mg/L"

# Row 388: E coli (>260)
$ws.Cells.Item(388, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(388, 2).Value = "E coli (>260)"
$ws.Cells.Item(388, 3).Value = "D"
$ws.Cells.Item(388, 4).Value = "2019 - 2023"
$ws.Cells.Item(388, 5).Value = "RepSite"
$ws.Cells.Item(388, 6).Value = 190
$ws.Cells.Item(388, 7).Value = 845.280701754386
$ws.Cells.Item(388, 8).Value = 9700
$ws.Cells.Item(388, 9).Value = 4382
$ws.Cells.Item(388, 10).Value = 31.5789473684211
$ws.Cells.Item(388, 11).Value = 42.1052631578947
$ws.Cells.Item(388, 12).Value = 68.5
$ws.Cells.Item(388, 13).Value = 1381
$ws.Cells.Item(388, 14).Value = 2864
$ws.Cells.Item(388, 15).Value = 1792183
$ws.Cells.Item(388, 16).Value = 5512989
$ws.Cells.Item(388, 17).Value = "Horowhenua District"
$ws.Cells.Item(388, 18).Value = "Manawatū"
$ws.Cells.Item(388, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(388, 20).Value = "Mana_13a"
$ws.Cells.Item(388, 21).Value = "% exceedances over 260/100 mL"

# Row 389: E coli (>540)
$ws.Cells.Item(389, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(389, 2).Value = "E coli (>540)"
$ws.Cells.Item(389, 3).Value = "E"
$ws.Cells.Item(389, 4).Value = "2019 - 2023"
$ws.Cells.Item(389, 5).Value = "RepSite"
$ws.Cells.Item(389, 6).Value = 190
$ws.Cells.Item(389, 7).Value = 845.280701754386
$ws.Cells.Item(389, 8).Value = 9700
$ws.Cells.Item(389, 9).Value = 4382
$ws.Cells.Item(389, 10).Value = 31.5789473684211
$ws.Cells.Item(389, 11).Value = 42.1052631578947
$ws.Cells.Item(389, 12).Value = 68.5
$ws.Cells.Item(389, 13).Value = 1381
$ws.Cells.Item(389, 14).Value = 2864
$ws.Cells.Item(389, 15).Value = 1792183
$ws.Cells.Item(389, 16).Value = 5512989
$ws.Cells.Item(389, 17).Value = "Horowhenua District"
$ws.Cells.Item(389, 18).Value = "Manawatū"
$ws.Cells.Item(389, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(389, 20).Value = "Mana_13a"
$ws.Cells.Item(389, 21).Value = "% exceedances over 540/100 mL"

# Row 390: E coli (Median)
$ws.Cells.Item(390, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(390, 2).Value = "E coli (Median)"
$ws.Cells.Item(390, 3).Value = "D"
$ws.Cells.Item(390, 4).Value = "2019 - 2023"
$ws.Cells.Item(390, 5).Value = "RepSite"
$ws.Cells.Item(390, 6).Value = 190
$ws.Cells.Item(390, 7).Value = 845.280701754386
$ws.Cells.Item(390, 8).Value = 9700
$ws.Cells.Item(390, 9).Value = 4382
$ws.Cells.Item(390, 10).Value = 31.5789473684211
$ws.Cells.Item(390, 11).Value = 42.1052631578947
$ws.Cells.Item(390, 12).Value = 68.5
$ws.Cells.Item(390, 13).Value = 1381
$ws.Cells.Item(390, 14).Value = 2864
$ws.Cells.Item(390, 15).Value = 1792183
$ws.Cells.Item(390, 16).Value = 5512989
$ws.Cells.Item(390, 17).Value = "Horowhenua District"
$ws.Cells.Item(390, 18).Value = "Manawatū"
$ws.Cells.Item(390, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(390, 20).Value = "Mana_13a"
$ws.Cells.Item(390, 21).Value = "E. coli/100 mL"

# Row 391: E coli (95th Percentile)
$ws.Cells.Item(391, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(391, 2).Value = "E coli (95th Percentile)"
$ws.Cells.Item(391, 3).Value = "E"
$ws.Cells.Item(391, 4).Value = "2019 - 2023"
$ws.Cells.Item(391, 5).Value = "RepSite"
$ws.Cells.Item(391, 6).Value = 190
$ws.Cells.Item(391, 7).Value = 845.280701754386
$ws.Cells.Item(391, 8).Value = 9700
$ws.Cells.Item(391, 9).Value = 4382
$ws.Cells.Item(391, 10).Value = 31.5789473684211
$ws.Cells.Item(391, 11).Value = 42.1052631578947
$ws.Cells.Item(391, 12).Value = 68.5
$ws.Cells.Item(391, 13).Value = 1381
$ws.Cells.Item(391, 14).Value = 2864
$ws.Cells.Item(391, 15).Value = 1792183
$ws.Cells.Item(391, 16).Value = 5512989
$ws.Cells.Item(391, 17).Value = "Horowhenua District"
$ws.Cells.Item(391, 18).Value = "Manawatū"
$ws.Cells.Item(391, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(391, 20).Value = "Mana_13a"
$ws.Cells.Item(391, 21).Value = "E. coli/100 mL"

# Row 392: Ammoniacal-N (95th Percentile)
$ws.Cells.Item(392, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(392, 2).Value = "Ammoniacal-N (95th Percentile)"
$ws.Cells.Item(392, 3).Value = "B"
$ws.Cells.Item(392, 4).Value = "2019 - 2023"
$ws.Cells.Item(392, 5).Value = "RepSite"
$ws.Cells.Item(392, 6).Value = 0.02755
$ws.Cells.Item(392, 7).Value = 0.0330503463097208
$ws.Cells.Item(392, 8).Value = 0.1157346472659
$ws.Cells.Item(392, 9).Value = 0.07671
$ws.Cells.Item(392, 10).Value = $null
$ws.Cells.Item(392, 11).Value = $null
$ws.Cells.Item(392, 12).Value = 0.02304
$ws.Cells.Item(392, 13).Value = 0.04836
$ws.Cells.Item(392, 14).Value = 0.06822
$ws.Cells.Item(392, 15).Value = 1792183
$ws.Cells.Item(392, 16).Value = 5512989
$ws.Cells.Item(392, 17).Value = "Horowhenua District"
$ws.Cells.Item(392, 18).Value = "Manawatū"
$ws.Cells.Item(392, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(392, 20).Value = "Mana_13a"
$ws.Cells.Item(392, 21).Value = "mg NH4-N/L"

# Row 393: Ammoniacal-N (Median)
$ws.Cells.Item(393, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(393, 2).Value = "Ammoniacal-N (Median)"
$ws.Cells.Item(393, 3).Value = "A"
$ws.Cells.Item(393, 4).Value = "2019 - 2023"
$ws.Cells.Item(393, 5).Value = "RepSite"
$ws.Cells.Item(393, 6).Value = 0.02755
$ws.Cells.Item(393, 7).Value = 0.0330503463097208
$ws.Cells.Item(393, 8).Value = 0.1157346472659
$ws.Cells.Item(393, 9).Value = 0.07671
$ws.Cells.Item(393, 10).Value = $null
$ws.Cells.Item(393, 11).Value = $null
$ws.Cells.Item(393, 12).Value = 0.02304
$ws.Cells.Item(393, 13).Value = 0.04836
$ws.Cells.Item(393, 14).Value = 0.06822
$ws.Cells.Item(393, 15).Value = 1792183
$ws.Cells.Item(393, 16).Value = 5512989
$ws.Cells.Item(393, 17).Value = "Horowhenua District"
$ws.Cells.Item(393, 18).Value = "Manawatū"
$ws.Cells.Item(393, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(393, 20).Value = "Mana_13a"
$ws.Cells.Item(393, 21).Value = "mg NH4-N/L"

# Row 394: Nitrate-N (95th Percentile)
$ws.Cells.Item(394, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(394, 2).Value = "Nitrate-N (95th Percentile)"
$ws.Cells.Item(394, 3).Value = "A"
$ws.Cells.Item(394, 4).Value = "2019 - 2023"
$ws.Cells.Item(394, 5).Value = "RepSite"
$ws.Cells.Item(394, 6).Value = 0.496
$ws.Cells.Item(394, 7).Value = 0.562614035087719
$ws.Cells.Item(394, 8).Value = 1.51
$ws.Cells.Item(394, 9).Value = 1.252
$ws.Cells.Item(394, 10).Value = $null
$ws.Cells.Item(394, 11).Value = $null
$ws.Cells.Item(394, 12).Value = 0.3115
$ws.Cells.Item(394, 13).Value = 0.87578
$ws.Cells.Item(394, 14).Value = 1.1488
$ws.Cells.Item(394, 15).Value = 1792183
$ws.Cells.Item(394, 16).Value = 5512989
$ws.Cells.Item(394, 17).Value = "Horowhenua District"
$ws.Cells.Item(394, 18).Value = "Manawatū"
$ws.Cells.Item(394, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(394, 20).Value = "Mana_13a"
$ws.Cells.Item(394, 21).Value = "mg NO3-N/L"

# Row 395: Nitrate-N (Median)
$ws.Cells.Item(395, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(395, 2).Value = "Nitrate-N (Median)"
$ws.Cells.Item(395, 3).Value = "A"
$ws.Cells.Item(395, 4).Value = "2019 - 2023"
$ws.Cells.Item(395, 5).Value = "RepSite"
$ws.Cells.Item(395, 6).Value = 0.496
$ws.Cells.Item(395, 7).Value = 0.562614035087719
$ws.Cells.Item(395, 8).Value = 1.51
$ws.Cells.Item(395, 9).Value = 1.252
$ws.Cells.Item(395, 10).Value = $null
$ws.Cells.Item(395, 11).Value = $null
$ws.Cells.Item(395, 12).Value = 0.3115
$ws.Cells.Item(395, 13).Value = 0.87578
$ws.Cells.Item(395, 14).Value = 1.1488
$ws.Cells.Item(395, 15).Value = 1792183
$ws.Cells.Item(395, 16).Value = 5512989
$ws.Cells.Item(395, 17).Value = "Horowhenua District"
$ws.Cells.Item(395, 18).Value = "Manawatū"
$ws.Cells.Item(395, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(395, 20).Value = "Mana_13a"
$ws.Cells.Item(395, 21).Value = "mg NO3-N/L"

# Row 396: Soluble Inorganic Nitrogen (95th Percentile)
$ws.Cells.Item(396, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(396, 2).Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Cells.Item(396, 3).Value = $null
$ws.Cells.Item(396, 4).Value = "2019 - 2023"
$ws.Cells.Item(396, 5).Value = "RepSite"
$ws.Cells.Item(396, 6).Value = 0.5685
$ws.Cells.Item(396, 7).Value = 0.642758620689655
$ws.Cells.Item(396, 8).Value = 1.57
$ws.Cells.Item(396, 9).Value = 1.3192
$ws.Cells.Item(396, 10).Value = $null
$ws.Cells.Item(396, 11).Value = $null
$ws.Cells.Item(396, 12).Value = 0.4105
$ws.Cells.Item(396, 13).Value = 0.96912
$ws.Cells.Item(396, 14).Value = 1.19874
$ws.Cells.Item(396, 15).Value = 1792183
$ws.Cells.Item(396, 16).Value = 5512989
$ws.Cells.Item(396, 17).Value = "Horowhenua District"
$ws.Cells.Item(396, 18).Value = "Manawatū"
$ws.Cells.Item(396, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(396, 20).Value = "Mana_13a"
$ws.Cells.Item(396, 21).Value = "g/m3"

# Row 397: Soluble Inorganic Nitrogen (Median)
$ws.Cells.Item(397, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(397, 2).Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Cells.Item(397, 3).Value = $null
$ws.Cells.Item(397, 4).Value = "2019 - 2023"
$ws.Cells.Item(397, 5).Value = "RepSite"
$ws.Cells.Item(397, 6).Value = 0.5685
$ws.Cells.Item(397, 7).Value = 0.642758620689655
$ws.Cells.Item(397, 8).Value = 1.57
$ws.Cells.Item(397, 9).Value = 1.3192
$ws.Cells.Item(397, 10).Value = $null
$ws.Cells.Item(397, 11).Value = $null
$ws.Cells.Item(397, 12).Value = 0.4105
$ws.Cells.Item(397, 13).Value = 0.96912
$ws.Cells.Item(397, 14).Value = 1.19874
$ws.Cells.Item(397, 15).Value = 1792183
$ws.Cells.Item(397, 16).Value = 5512989
$ws.Cells.Item(397, 17).Value = "Horowhenua District"
$ws.Cells.Item(397, 18).Value = "Manawatū"
$ws.Cells.Item(397, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(397, 20).Value = "Mana_13a"
$ws.Cells.Item(397, 21).Value = "g/m3"

# Row 398: Total Nitrogen (95th Percentile)
$ws.Cells.Item(398, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(398, 2).Value = "Total Nitrogen (95th Percentile)"
$ws.Cells.Item(398, 3).Value = $null
$ws.Cells.Item(398, 4).Value = "2019 - 2023"
$ws.Cells.Item(398, 5).Value = "RepSite"
$ws.Cells.Item(398, 6).Value = 0.76
$ws.Cells.Item(398, 7).Value = 0.893684210526316
$ws.Cells.Item(398, 8).Value = 1.83
$ws.Cells.Item(398, 9).Value = 1.76
$ws.Cells.Item(398, 10).Value = $null
$ws.Cells.Item(398, 11).Value = $null
$ws.Cells.Item(398, 12).Value = 0.615
$ws.Cells.Item(398, 13).Value = 1.2143
$ws.Cells.Item(398, 14).Value = 1.554
$ws.Cells.Item(398, 15).Value = 1792183
$ws.Cells.Item(398, 16).Value = 5512989
$ws.Cells.Item(398, 17).Value = "Horowhenua District"
$ws.Cells.Item(398, 18).Value = "Manawatū"
$ws.Cells.Item(398, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(398, 20).Value = "Mana_13a"
$ws.Cells.Item(398, 21).Value = "g/m3"

# Row 399: Total Nitrogen (Median)
$ws.Cells.Item(399, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(399, 2).Value = "Total Nitrogen (Median)"
$ws.Cells.Item(399, 3).Value = $null
$ws.Cells.Item(399, 4).Value = "2019 - 2023"
$ws.Cells.Item(399, 5).Value = "RepSite"
$ws.Cells.Item(399, 6).Value = 0.76
$ws.Cells.Item(399, 7).Value = 0.893684210526316
$ws.Cells.Item(399, 8).Value = 1.83
$ws.Cells.Item(399, 9).Value = 1.76
$ws.Cells.Item(399, 10).Value = $null
$ws.Cells.Item(399, 11).Value = $null
$ws.Cells.Item(399, 12).Value = 0.615
$ws.Cells.Item(399, 13).Value = 1.2143
$ws.Cells.Item(399, 14).Value = 1.554
$ws.Cells.Item(399, 15).Value = 1792183
$ws.Cells.Item(399, 16).Value = 5512989
$ws.Cells.Item(399, 17).Value = "Horowhenua District"
$ws.Cells.Item(399, 18).Value = "Manawatū"
$ws.Cells.Item(399, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(399, 20).Value = "Mana_13a"
$ws.Cells.Item(399, 21).Value = "g/m3"

# Row 400: Total Phosphorus (95th Percentile)
$ws.Cells.Item(400, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(400, 2).Value = "Total Phosphorus (95th Percentile)"
$ws.Cells.Item(400, 3).Value = $null
$ws.Cells.Item(400, 4).Value = "2019 - 2023"
$ws.Cells.Item(400, 5).Value = "RepSite"
$ws.Cells.Item(400, 6).Value = 0.064
$ws.Cells.Item(400, 7).Value = 0.0938245614035088
$ws.Cells.Item(400, 8).Value = 0.362
$ws.Cells.Item(400, 9).Value = 0.298
$ws.Cells.Item(400, 10).Value = $null
$ws.Cells.Item(400, 11).Value = $null
$ws.Cells.Item(400, 12).Value = 0.0645
$ws.Cells.Item(400, 13).Value = 0.13248
$ws.Cells.Item(400, 14).Value = 0.20782
$ws.Cells.Item(400, 15).Value = 1792183
$ws.Cells.Item(400, 16).Value = 5512989
$ws.Cells.Item(400, 17).Value = "Horowhenua District"
$ws.Cells.Item(400, 18).Value = "Manawatū"
$ws.Cells.Item(400, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(400, 20).Value = "Mana_13a"
$ws.Cells.Item(400, 21).Value = "g/m3"

# Row 401: Total Phosphorus (Median)
$ws.Cells.Item(401, 1).Value = "Manawatu at Whirokino"
$ws.Cells.Item(401, 2).Value = "Total Phosphorus (Median)"
$ws.Cells.Item(401, 3).Value = $null
$ws.Cells.Item(401, 4).Value = "2019 - 2023"
$ws.Cells.Item(401, 5).Value = "RepSite"
$ws.Cells.Item(401, 6).Value = 0.064
$ws.Cells.Item(401, 7).Value = 0.0938245614035088
$ws.Cells.Item(401, 8).Value = 0.362
$ws.Cells.Item(401, 9).Value = 0.298
$ws.Cells.Item(401, 10).Value = $null
$ws.Cells.Item(401, 11).Value = $null
$ws.Cells.Item(401, 12).Value = 0.0645
$ws.Cells.Item(401, 13).Value = 0.13248
$ws.Cells.Item(401, 14).Value = 0.20782
$ws.Cells.Item(401, 15).Value = 1792183
$ws.Cells.Item(401, 16).Value = 5512989
$ws.Cells.Item(401, 17).Value = "Horowhenua District"
$ws.Cells.Item(401, 18).Value = "Manawatū"
$ws.Cells.Item(401, 19).Value = "Coastal Manawatu"
$ws.Cells.Item(401, 20).Value = "Mana_13a"
$ws.Cells.Item(401, 21).Value = "g/m3"

